$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> new value for column F (dSF)
$updates = @{
    2  = -2
    3  = -5
    4  = -2
    5  = 0
    6  = -3
    7  = -1
    8  = 4
    9  = -6
    10 = 7
    11 = -2
    12 = -2
    13 = -1
    15 = 7
    16 = -4
    17 = 2
    18 = -1
    19 = 7
    20 = -4
    22 = -1
    23 = 1
    24 = 6
    25 = 3
    26 = -3
    27 = -2
    28 = -2
    29 = -1
    31 = 6
    32 = 3
    33 = 1
    34 = -3
    35 = 3
    36 = 0
    37 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
